$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Translate header row (Portuguese -> English) ---
$ws.Range("A1").Value = "email"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "value"
$ws.Range("D1").Value = "product"
$ws.Range("E1").Value = "amount"
$ws.Range("F1").Value = "category"

# --- Translate data rows (rows 2-21) ---
$letters = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
$categoryCycle = @("category1","category2","category3")

for ($i = 0; $i -lt 20; $i++) {
    $row = $i + 2

    # Column A: email (row 2 stays blank, same as before)
    if ($row -gt 2) {
        $n = $row - 1
        $ws.Cells.Item($row, 1).Value = "user" + $n + "@example.com"
    }

    # Column D: product name
    $ws.Cells.Item($row, 4).Value = "Product " + $letters[$i]

    # Column F: category
    $ws.Cells.Item($row, 6).Value = $categoryCycle[$i % 3]
}

# --- Fix the sign on C8 (-300 -> 300) ---
$ws.Range("C8").Value = 300

# --- Clear the product text on D15, keeping it as a quote-prefixed empty cell ---
$ws.Range("D15").Value = "'"

# --- Shrink the data font from 11pt to 10pt across the whole table ---
$ws.Range("A1:F21").Font.Size = 10
